$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-12-22 Monday" "2025-12-23 Tuesday"

Replace-Text "38×16=608" "62×43=2666"
Replace-Text "55×91=5005" "76×76=5776"
Replace-Text "24×61=1464" "24×88=2112"
Replace-Text "65×45=2925" "99×98=9702"
Replace-Text "62×82=5084" "66×30=1980"

Replace-Text "37×88=3256" "98×45=4410"
Replace-Text "27×84=2268" "20×53=1060"
Replace-Text "38×72=2736" "74×35=2590"
Replace-Text "71×80=5680" "11×79=869"
Replace-Text "16×53=848" "15×32=480"

Replace-Text "69×23=1587" "85×65=5525"
Replace-Text "91×30=2730" "27×56=1512"
Replace-Text "64×68=4352" "93×96=8928"
Replace-Text "48×18=864" "50×91=4550"
Replace-Text "72×58=4176" "39×70=2730"

Replace-Text "85×35=2975" "97×86=8342"
Replace-Text "22×19=418" "58×51=2958"
Replace-Text "36×50=1800" "64×90=5760"
Replace-Text "96×63=6048" "39×70=2730"
Replace-Text "61×68=4148" "85×19=1615"

Replace-Text "60×95=5700" "28×60=1680"
Replace-Text "91×56=5096" "49×96=4704"
Replace-Text "60×70=4200" "62×40=2480"
Replace-Text "71×14=994" "34×27=918"
Replace-Text "35×56=1960" "48×53=2544"
